# Edit script: rename balance tracking columns and refresh chain balance data
# C: USDT -> ETH_arb (reuses prior ETH_arb column E data)
# D: USDC -> ETH_linea (new data)
# E: ETH_arb -> ETH_op (new data)
# F: USDT_arb -> ETH_zksync (new data, replaces USDT error messages)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---
$ws.Range("C1").Value = "ETH_arb"
$ws.Range("D1").Value = "ETH_linea"
$ws.Range("E1").Value = "ETH_op"
$ws.Range("F1").Value = "ETH_zksync"

# --- Per-row new data for columns D (ETH_linea), E (ETH_op), F (ETH_zksync) ---
# Each entry: row, newD, newE, newF (newD = "0" means leave existing numeric 0 as-is)
$rowData = @(
    @(2, "0", "0.002200836190495412", "0.015099496454615183"),
    @(3, "0", "0.00521524178559823", "0.015176544130623508"),
    @(4, "0", "0.002135100732171918", "0.015176489122341696"),
    @(5, "0", "0.001624884842310747", "0.015172670325214127"),
    @(6, "0", "0.001895215160174945", "0.015125857492476492"),
    @(7, "0", "0.001488370167283218", "0.015149127631979354"),
    @(8, "0", "0.002053013326231324", "0.015116798681313845"),
    @(9, "0", "0.001983289096465981", "0.015198702062870001"),
    @(10, "0", "0.001663123103873953", "0.015255375021637252"),
    @(11, "0", "0.001483688916342881", "0.015172954782635143"),
    @(12, "0", "0.001652235479665607", "0.015180848922078652"),
    @(13, "0", "0.002137748200075364", "0.015160322659445499"),
    @(14, "0", "0.001958175685403652", "0.015103316942946463"),
    @(15, "0", "0.001296491744178415", "0.016433254939550414"),
    @(16, "0", "0.005556224387734968", "0.015142242722078457"),
    @(17, "0", "0.006771998882783559", "0.015193494536207619"),
    @(18, "0", "0.00267645163747792", "0.015128569538923159"),
    @(19, "0", "0.001864491364989403", "0.015124818183854209"),
    @(20, "0", "0.001985119617351436", "0.015184096214549963"),
    @(21, "0", "0.001932794729509414", "0.015185012825205172"),
    @(22, "0", "0.001534215401876905", "0.015137931358804504"),
    @(23, "0", "0.006988865718991985", "0.015097292909320977"),
    @(24, "0", "0.006965630793153929", "0.015081014434191612"),
    @(25, "0", "0.006946357939342155", "0.015146750449299984"),
    @(26, "0", "0.001585244504725384", "0.015061511737303648"),
    @(27, "0", "0.006976246037370431", "0.015128902268251703"),
    @(28, "0", "0.028968675830108055", "0.015085647067183624"),
    @(29, "0", "0.005367869708379562", "0.015281110689133648"),
    @(30, "0", "0.001755738659213979", "0.015189404034243648"),
    @(31, "0", "0.006896193222684303", "0.015090874380522679"),
    @(32, "0", "0.006535266711161772", "0.015070396443533648"),
    @(33, "0", "0.006155333459011032", "0.015213845002851701"),
    @(34, "0", "0.006906871344373532", "0.01511098021385267"),
    @(35, "0", "0.001514007941568989", "0.01504176328428265"),
    @(36, "0", "0.001545083728792859", "0.0151344382245207"),
    @(37, "0", "0.001326864850068064", "0.01514273026063265"),
    @(38, "0", "0.001725108721096234", "0.016338787896783631"),
    @(39, "0", "0.001695754682153717", "0.015148099980023624"),
    @(40, "0", "0.001728408079407527", "0.015086441825633624"),
    @(41, "0", "0.001685733263214249", "0.015110808620653624"),
    @(42, "0", "0.000686693733995163", "0.01515590125"),
    @(43, "0", "0.001434460115268412", "0.015186650543741681"),
    @(44, "0", "0.001505239975768865", "0.015082961982542655"),
    @(45, "0", "0.001396992005637018", "0.015184672914401679"),
    @(46, "0", "0.001365037272692224", "0.015179290907403624"),
    @(47, "0", "0.001377075770320885", "0.015103386821299719"),
    @(48, "0", "0.001336206939379119", "0.015103875136222641"),
    @(49, "0", "0.001323234685152189", "0.015241652762603622"),
    @(50, "0", "0.001444960922498667", "0.015231633624792641"),
    @(51, "0", "0.00131202094708996", "0.0150901723359507"),
    @(52, "0.02846647051653649", "0.001473325559169078", "0.015165315698868709"),
    @(53, "0.028425680214173871", "0.001352156970233569", "0.01520021724430667"),
    @(54, "0.028237495696581474", "0.001080181904172627", "0.015130150280226875"),
    @(55, "0.028419199010261388", "0.00146239008298656", "0.015103229927422268"),
    @(56, "0.028261909648041693", "0.001274187714531556", "0.015234557284383141"),
    @(57, "0.028388327285206513", "0.001796139665146785", "0.015172247841539798"),
    @(58, "0.028450647802057852", "0.001273762886168379", "0.01511744237770285"),
    @(59, "0.028126773102295135", "0.001708602401400727", "0.015180773999533216"),
    @(60, "0.028401624409163199", "0.001109036269904536", "0.015195073610463794"),
    @(61, "0.028333888194507543", "0.001246554603928167", "0.015207498875896864"),
    @(62, "0.02814781025045381", "0.001259801193973763", "0.000085661225003894"),
    @(63, "0.028328089512886782", "0.001059335078061579", "0.015126611999426629"),
    @(64, "0.028093687292859665", "0.001332055232610982", "0.015147526015047197"),
    @(65, "0.02821669539371679", "0.00115281824477365", "0.015197011862548334"),
    @(66, "0.028147815413161124", "0.001285881855635335", "0.015217209738616216"),
    @(67, "0.028473859060892187", "0.001833304936282787", "0.015087071536103347"),
    @(68, "0.027325186303633706", "0.001255004833943484", "0.015204269453907454"),
    @(69, "0.032721856309667415", "0.001216295913372081", "0.015085888476090345"),
    @(70, "0.032168055438447103", "0.001933206533423739", "0.0151597194417841"),
    @(71, "0.030261890882806431", "0.00182061041733268", "0.015179520508399497"),
    @(72, "0.028617022712428566", "0.001489613072390388", "0.015150373058915061"),
    @(73, "0.02952097033314789", "0.001511124285991686", "0.015190575781436488"),
    @(74, "0.02883347197225529", "0.001146886792482124", "0.015152643773083303"),
    @(75, "0.028542545315584325", "0.001138283043640316", "0.015174670423348568"),
    @(76, "0.027874395111667421", "0.001686374953623174", "0.015198119460064174"),
    @(77, "0.027960244500396393", "0.001090854388742621", "0.015127381205620652"),
    @(78, "0.028057593556592169", "0.001212317632533659", "0.015107748702400526"),
    @(79, "0.028474571932921808", "0.001227610490510062", "0.015192932611938536"),
    @(80, "0.02713182051301753", "0.001436726010288015", "0.015210373002408372"),
    @(81, "0.027365933878516829", "0.001027979331134178", "0.015219242318080282"),
    @(82, "0.027332330859437399", "0.001912749857251252", "0.015096423074437814"),
    @(83, "0.02715017147482356", "0.001053473250602096", "0.015149234939587917"),
    @(84, "0.026942297626387895", "0.001325625319249775", "0.015131780477591974"),
    @(85, "0.028490045552660349", "0.001351362026756053", "0.015219009764980229"),
    @(86, "0.02942061163559607", "0.001073624958715281", "0.0151152049716855"),
    @(87, "0.027487614117409531", "0.00103322510816264", "0.015195177693997175"),
    @(88, "0.028120142638262229", "0.001027417379635377", "0.015090989781468095"),
    @(89, "0.027724594735270515", "0.001049907140686047", "0.015198081567062105"),
    @(90, "0.028618121624058864", "0.001104711605733498", "0.015133512050852379"),
    @(91, "0.028488359043567372", "0.001722428409873749", "0.015212502995313673"),
    @(92, "0", "0.001074563621035278", "0.015164729055557709"),
    @(93, "6.3E-14", "0.00104892441152464", "0.015209413303378703"),
    @(94, "0.000005111190961017", "0.001594045073531235", "0.015222546484549671"),
    @(95, "0.028528528278438175", "0.001479136022820603", "0.015137512455766257"),
    @(96, "0.02798139187608166", "0.000432361768235113", "0.015122223135242137"),
    @(97, "0.000004427723890012", "0.001768188802595249", "0.01517650334745408"),
    @(98, "0.000019241062219013", "0.00119823435570017", "0.015217738936845375"),
    @(99, "0.000014949919542012", "0.001604026741577399", "0.015111291330032973"),
    @(100, "0.02770965844091444", "0.001229670788449995", "0.015232834678531705"),
    @(101, "0.000005766893349011", "0.000269457736299954", "0.0151891052910341")
)

foreach ($entry in $rowData) {
    $r = $entry[0]
    $newD = $entry[1]
    $newE = $entry[2]
    $newF = $entry[3]

    # Column C (ETH_arb) takes over the value that used to live in column E (old ETH_arb)
    $oldE = $ws.Cells.Item($r, 5).Value()
    $cCell = $ws.Cells.Item($r, 3)
    $cCell.NumberFormat = "@"
    $cCell.Value = $oldE
    $cCell.ClearFormats()

    # Column D (ETH_linea) - only rewrite when the diff actually changes it
    if ($newD -ne "0") {
        $dCell = $ws.Cells.Item($r, 4)
        $dCell.NumberFormat = "@"
        $dCell.Value = $newD
        $dCell.ClearFormats()
    }

    # Column E (ETH_op)
    $eCell = $ws.Cells.Item($r, 5)
    $eCell.NumberFormat = "@"
    $eCell.Value = $newE
    $eCell.ClearFormats()

    # Column F (ETH_zksync)
    $fCell = $ws.Cells.Item($r, 6)
    $fCell.NumberFormat = "@"
    $fCell.Value = $newF
    $fCell.ClearFormats()
}
